# Scheduled market-data refresh for the Leve profit sheets.
# Updates currentAveragePrice(NQ/HQ) (H:J), LevePrice(NQ/HQ) (K:L) and the
# derived LeveProfit(NQ/HQ) (M:N) columns with freshly polled prices; a few
# rows gain/lose M/N cells entirely depending on whether NQ/HQ listings exist.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 42: Eye of the Beholder / Hi-Potion of Dexterity
$ws.Range("H42").Value = 86.75
$ws.Range("I42").Value = 79.333336
$ws.Range("K42").Value = 238.000008
$ws.Range("M42").Value = -8.000008000000008
# Row 61: Not Taking No for an Answer / Mega-Potion of Strength
$ws.Range("H61").Value = 608.8333
$ws.Range("I61").Value = 608.8333
$ws.Range("K61").Value = 1826.4999
$ws.Range("M61").Value = -1654.4999
# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 113180.336
$ws.Range("I76").Value = 113180.336
$ws.Range("K76").Value = 113180.336
$ws.Range("M76").Value = -112865.336
# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 113180.336
$ws.Range("I79").Value = 113180.336
$ws.Range("K79").Value = 113180.336
$ws.Range("M79").Value = -112088.336
# Row 87: There Was a Late Fee / Noble Gold
$ws.Range("H87").Value = 54999.5
$ws.Range("J87").Value = 54999.5
$ws.Range("L87").Value = 54999.5
$ws.Range("N87").Value = -57495.5
# Row 90: A Gate Arcane Is Dragon's Bane (L) / Noble Gold
$ws.Range("H90").Value = 54999.5
$ws.Range("J90").Value = 54999.5
$ws.Range("L90").Value = 164998.5
$ws.Range("N90").Value = -177478.5
# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 7156.706
$ws.Range("I98").Value = 6979.625
$ws.Range("K98").Value = 6979.625
$ws.Range("M98").Value = -5481.625
# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 40747330
$ws.Range("I113").Value = 18520854
$ws.Range("J113").Value = 55564980
$ws.Range("K113").Value = 18520854
$ws.Range("L113").Value = 55564980
$ws.Range("M113").Value = -18517600
$ws.Range("N113").Value = -55571488
# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 7156.706
$ws.Range("I122").Value = 6979.625
$ws.Range("K122").Value = 20938.875
$ws.Range("M122").Value = -18488.875
# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 83336230
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 127: Liquid Competence / Competent Craftsman's Draught
$ws.Range("H127").Value = 3320.3572
$ws.Range("I127").Value = 2207.0833
$ws.Range("K127").Value = 6621.249899999999
$ws.Range("M127").Value = -1661.249899999999
# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 1661.125
$ws.Range("I141").Value = 1727.1428
$ws.Range("K141").Value = 5181.428400000001
$ws.Range("M141").Value = -1.428400000000693

$ws = $wb.Worksheets.Item("ARM")
# Row 28: 246 Kinds of Cheese / Iron Frypan
$ws.Range("H28").Value = 14749.2
$ws.Range("I28").Value = 4951
$ws.Range("K28").Value = 4951
$ws.Range("M28").Value = -4759
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 1925684.4
$ws.Range("I32").Value = 2051658.8
$ws.Range("K32").Value = 2051658.8
$ws.Range("M32").Value = -2051371.8
# Row 41: Skillet Scandal / White Skillet
$ws.Range("H41").Value = 632
$ws.Range("I41").Value = 632
$ws.Range("K41").Value = 632
$ws.Range("M41").Value = -218
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 7452.154
$ws.Range("J45").Value = 11839.571
$ws.Range("L45").Value = 11839.571
$ws.Range("N45").Value = -12593.571
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 5481.2554
$ws.Range("I61").Value = 2437.027
$ws.Range("K61").Value = 2437.027
$ws.Range("M61").Value = -2225.027
# Row 99: Home Cooking / Doman Iron Frypan
$ws.Range("H99").Value = 14749.2
$ws.Range("I99").Value = 4951
$ws.Range("K99").Value = 4951
$ws.Range("M99").Value = -1956
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 7594.2646
$ws.Range("I132").Value = 6366.364
$ws.Range("K132").Value = 19099.092
$ws.Range("M132").Value = -16569.092
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 5481.2554
$ws.Range("I136").Value = 2437.027
$ws.Range("K136").Value = 7311.081
$ws.Range("M136").Value = -4761.081
# Row 138: Don't Ask about the Rivets / Titanium Gold Helm of Casting
$ws.Range("H138").Value = 79807.664
$ws.Range("J138").Value = 79884.39999999999
$ws.Range("L138").Value = 79884.39999999999
$ws.Range("N138").Value = -90164.39999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 24: Honest Ballast / Initiate's Head Knife
$ws.Range("H24").Value = 500
$ws.Range("I24").Value = 500
$ws.Range("K24").Value = 500
$ws.Range("M24").Value = -265
# Row 29: Powderpost Derby / Initiate's Saw
$ws.Range("H29").Value = 950
$ws.Range("I29").Value = 900
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 900
$ws.Range("L29").Value = 1000
$ws.Range("M29").Value = -611
$ws.Range("N29").Value = -1578
# Row 36: I Saw What You Did There / Iron Chocobotail Saw
$ws.Range("H36").Value = 393.5
$ws.Range("I36").Value = 393.5
$ws.Range("K36").Value = 393.5
$ws.Range("M36").Value = 140.5
# Row 54: Get Me to the War on Time / Cobalt Doming Hammer
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
# Row 96: Hammer Time / High Steel Sledgehammer
$ws.Range("H96").Value = 30453.625
$ws.Range("I96").Value = 16770.4
$ws.Range("J96").Value = 53259
$ws.Range("K96").Value = 16770.4
$ws.Range("L96").Value = 53259
$ws.Range("M96").Value = -14024.4
$ws.Range("N96").Value = -58751
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 1778.6296
$ws.Range("I105").Value = 1240.8334
$ws.Range("K105").Value = 1240.8334
$ws.Range("M105").Value = 506.1666

$ws = $wb.Worksheets.Item("CRP")
# Row 55: Ready for a Rematch / Mythril Lance
$ws.Range("H55").Value = 23081
$ws.Range("J55").Value = 23081
$ws.Range("L55").Value = 23081
$ws.Range("N55").Value = -23711
# Row 93: Reeling for Rods / Muudhorn Fishing Rod
$ws.Range("H93").Value = 32235.5
$ws.Range("I93").Value = 25000
$ws.Range("J93").Value = 53942
$ws.Range("K93").Value = 25000
$ws.Range("L93").Value = 53942
$ws.Range("M93").Value = -23128
$ws.Range("N93").Value = -57686
# Row 141: No Greater Treasure / Claro Walnut Necklace of Gathering
$ws.Range("H141").Value = 60827.2
$ws.Range("J141").Value = 60827.2
$ws.Range("L141").Value = 60827.2
$ws.Range("N141").Value = -71187.2

$ws = $wb.Worksheets.Item("CUL")
# Row 80: Saucy for a Suitor / Hollandaise Sauce
$ws.Range("H80").Value = 200003300
$ws.Range("I80").Value = 250002000
$ws.Range("J80").Value = 166670830
$ws.Range("K80").Value = 750006000
$ws.Range("L80").Value = 500012490
$ws.Range("M80").Value = -750005064
$ws.Range("N80").Value = -500014362
# Row 83: Saved by the Sauce (L) / Hollandaise Sauce
$ws.Range("H83").Value = 200003300
$ws.Range("I83").Value = 250002000
$ws.Range("J83").Value = 166670830
$ws.Range("K83").Value = 2250018000
$ws.Range("L83").Value = 1500037470
$ws.Range("M83").Value = -2250013320
$ws.Range("N83").Value = -1500046830

$ws = $wb.Worksheets.Item("GSM")
# Row 134: Guaranteed Gem / Ihuykanite
$ws.Range("H134").Value = 119326
$ws.Range("J134").Value = 119326
$ws.Range("L134").Value = 357978
$ws.Range("N134").Value = -363048

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 7433.1177
$ws.Range("I7").Value = 6635.4
$ws.Range("K7").Value = 6635.4
$ws.Range("M7").Value = -6523.4
# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 4430
$ws.Range("I68").Value = 3628.7144
$ws.Range("K68").Value = 3628.7144
$ws.Range("M68").Value = -2879.7144
# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 4430
$ws.Range("I71").Value = 3628.7144
$ws.Range("K71").Value = 18143.572
$ws.Range("M71").Value = -14399.572
# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 7433.1177
$ws.Range("I126").Value = 6635.4
$ws.Range("K126").Value = 19906.2
$ws.Range("M126").Value = -17436.2
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 20841848
$ws.Range("I132").Value = 41669860
$ws.Range("J132").Value = 13833.333
$ws.Range("K132").Value = 125009580
$ws.Range("L132").Value = 41499.999
$ws.Range("M132").Value = -125007050
$ws.Range("N132").Value = -46559.999

$ws = $wb.Worksheets.Item("WVR")
# Row 14: Hat in Hand / Straw Hat
$ws.Range("H14").Value = 333335360
$ws.Range("J14").Value = 4000
$ws.Range("L14").Value = 4000
$ws.Range("N14").Value = -4336
# Row 26: New Shoes, New Me / Cotton Dress Shoes
$ws.Range("H26").Value = 3000.5
$ws.Range("I26").Value = 1001
$ws.Range("J26").Value = 5000
$ws.Range("K26").Value = 1001
$ws.Range("L26").Value = 5000
$ws.Range("M26").Value = -708
$ws.Range("N26").Value = -5586
# Row 56: Full Moon Fever / Felt Chausses
$ws.Range("H56").Value = 36091.668
$ws.Range("I56").Value = 36637.5
$ws.Range("J56").Value = 35000
$ws.Range("K56").Value = 36637.5
$ws.Range("L56").Value = 35000
$ws.Range("M56").Value = -35923.5
$ws.Range("N56").Value = -36428
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 108951.695
$ws.Range("J122").Value = 5139
$ws.Range("L122").Value = 15417
$ws.Range("N122").Value = -20317
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 14290389
$ws.Range("I132").Value = 19235786
$ws.Range("K132").Value = 57707358
$ws.Range("M132").Value = -57704828
